$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values; C2 and E2 become empty (removed)
$ws.Range("B2").Value = 1.9341576716069004
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.7960150970981306
$ws.Range("E2").ClearContents()

# Update row 3 values
$ws.Range("B3").Value = 1.6688906381123014
$ws.Range("C3").Value = -1.3172335172622707
$ws.Range("D3").Value = 2.2270202844390838
$ws.Range("E3").Value = -2.0383289913388207

# Update selection to match new reduced range B1:E3
$ws.Range("B1:E3").Select()
